$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.214.04"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "1.682.60"
$ws.Range("E3").Value = "  +0.22%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "216.28"
$ws.Range("E5").Value = "  -0.76%  "

$ws.Range("D6").Value = "0.5279"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("E8").Value = "  +0.34%  "

$ws.Range("D9").Value = "0.06366"
$ws.Range("E9").Value = "  -1.71%  "

$ws.Range("E10").Value = "  -2.12%  "

$ws.Range("D11").Value = "0.07615"
$ws.Range("E11").Value = "  +0.95%  "

$ws.Range("D12").Value = "1.703.11"
$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("D13").Value = "4.524"
$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("D14").Value = "0.5762"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("D15").Value = "0.000008241"
$ws.Range("E15").Value = "  -2.57%  "

$ws.Range("D16").Value = "66.39"
$ws.Range("E16").Value = "  +2.46%  "

$ws.Range("D17").Value = "26.243.80"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("D19").Value = "4.868"
$ws.Range("E19").Value = "  -0.71%  "

$ws.Range("D20").Value = "10.76"
$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("D22").Value = "6.242"
$ws.Range("E22").Value = "  +0.59%  "

$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").Value = "149.19"
$ws.Range("E24").Value = "  +2.42%  "

$ws.Range("D25").Value = "0.1263"
$ws.Range("E25").Value = "  -1.22%  "

$ws.Range("D26").Value = "7.722"

$ws.Range("D27").Value = "15.86"
$ws.Range("E27").Value = "  +0.61%  "

$ws.Range("D28").Value = "0.06407"
$ws.Range("E28").Value = "  -1.33%  "

$ws.Range("D29").Value = "1.376"
$ws.Range("E29").Value = "  -0.52%  "

$ws.Range("E30").Value = "  -0.53%  "

$ws.Range("D31").Value = "3.568"
$ws.Range("E31").Value = "  -0.28%  "

$ws.Range("D32").Value = "3.567"
$ws.Range("E32").Value = "  -0.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.680"
$ws.Range("E33").Value = "  +0.77%  "

$ws.Range("E34").Value = "  -1.01%  "

$ws.Range("D35").Value = "0.6132"
$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("D36").Value = "2.417"
$ws.Range("E36").Value = "  +0.69%  "

$ws.Range("D37").Value = "2.749"
$ws.Range("E37").Value = "  +1.78%  "

$ws.Range("D38").Value = "6.178"
$ws.Range("E38").Value = "  -1.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01620"
$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("D40").Value = "1.096.03"
$ws.Range("E40").Value = "  -1.43%  "

$ws.Range("D41").Value = "0.8835"

$ws.Range("E42").Value = "  -0.45%  "

$ws.Range("D43").Value = "100.46"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").Value = "1.833.56"
$ws.Range("E44").Value = "  +0.26%  "

$ws.Range("D45").Value = "0.00000000111"
$ws.Range("E45").Value = "  +3.19%  "

$ws.Range("D46").Value = "57.49"
$ws.Range("E46").Value = "  +0.74%  "

$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("D48").Value = "8.079"
$ws.Range("E48").Value = "  -0.98%  "

$ws.Range("D49").Value = "0.05267"
$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("E50").Value = "  -0.27%  "

$ws.Range("D51").Value = "6.016"
$ws.Range("E51").Value = "  -0.97%  "
